$d = $word.ActiveDocument
$wordNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ------------------------------------------------------------------
# 1) Split the run "Nosso curso de GIT E GITHUB !!!" so that
#    "GITHUB !!!" is wrapped with gramStart/gramEnd proofErr markers,
#    matching a grammar-check artifact Word would normally add.
# ------------------------------------------------------------------
$findRange = $d.Content
$found = $findRange.Find.Execute("GITHUB !!!", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target text 'GITHUB !!!'"
}

# Re-wrap the hit in a fresh Range; InsertXML on a Range that is still
# "attached" to a Find operation does not replace the found text, it
# appends after it instead.
$targetRange = $d.Range($findRange.Start, $findRange.End)

$gramXml = '<w:p xmlns:w="' + $wordNs + '"><w:proofErr w:type="gramStart"/><w:r><w:t>GITHUB !!!</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>'
$targetRange.InsertXML($gramXml)

# ------------------------------------------------------------------
# 2) Add a new paragraph "Aprendemos Sobre git" right after the
#    paragraph that now ends in "...GITHUB !!!", moving the
#    "_GoBack" bookmark along so it stays at the very end of the
#    document content (as it was before the edit).
# ------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bmRange = $bm.Range
$para = $bmRange.Paragraphs(1)
$paraEnd = $para.Range.End

# Insert a new paragraph mark right after the current paragraph, then
# write the new paragraph's text plus a temporary two-character
# placeholder. The placeholder keeps the bookmark's insertion point
# away from the position immediately before a paragraph mark, which
# this runtime mishandles when adding a brand new bookmark there.
$insertionPoint = $d.Range($paraEnd - 1, $paraEnd - 1)
$insertionPoint.InsertParagraphAfter()

$newParaStart = $paraEnd
$newTextRange = $d.Range($newParaStart, $newParaStart)
$newText = "Aprendemos Sobre git"
$placeholder = "ZZ"
$newTextRange.InsertAfter($newText + $placeholder)

$bookmarkPos = $newParaStart + $newText.Length
$safeTarget = $d.Range($bookmarkPos, $bookmarkPos)

$d.Bookmarks.Item("_GoBack").Delete()
$d.Bookmarks.Add("_GoBack", $safeTarget)

$placeholderRange = $d.Range($bookmarkPos, $bookmarkPos + $placeholder.Length)
$placeholderRange.Delete()
